# Updated symbol list on Sun Dec 18 09:43:01 UTC 2022 with GitHub Actions
# Applies the per-row "Price" (column D) and a couple of "Volume(1h)" (column E)
# text updates captured by the diff. Source cells are stored as text
# (inline strings) even though many values look numeric, so every
# numeric-looking replacement is written with a leading apostrophe to force
# Excel to keep storing it as text rather than silently re-typing it as a
# number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") updates
$ws.Range("D2").Value  = "'245.47"
$ws.Range("D3").Value  = "'22.70"
$ws.Range("D4").Value  = "'5.533"
$ws.Range("D5").Value  = "'0.05605"
$ws.Range("D6").Value  = "'6.473"
$ws.Range("D7").Value  = "'0.8032"
$ws.Range("D8").Value  = "'1.058"
$ws.Range("D9").Value  = "'0.1421"
$ws.Range("D10").Value = "'0.07283"
$ws.Range("D11").Value = "'0.03201"
$ws.Range("D12").Value = "'0.02986"
$ws.Range("D13").Value = "'0.09260"
$ws.Range("D14").Value = "'0.001665"
$ws.Range("D15").Value = "'2.972"
$ws.Range("D16").Value = "'0.04701"
$ws.Range("D17").Value = "'0.0005991"
$ws.Range("D18").Value = "'0.006290"
$ws.Range("D19").Value = "'0.001056"
$ws.Range("D20").Value = "'0.003796"
$ws.Range("D21").Value = "'0.0001505"
$ws.Range("D22").Value = "'0.0004015"
$ws.Range("D24").Value = "'3.407"
$ws.Range("D27").Value = "'0.1294"
$ws.Range("D40").Value = "'0.04171"
$ws.Range("D41").Value = "'0.006981"
$ws.Range("D42").Value = "'0.1039"
$ws.Range("D43").Value = "'0.003090"
$ws.Range("D44").Value = "'0.009287"
$ws.Range("D45").Value = "'0.00005672"
$ws.Range("D46").Value = "'0.00000000753"
$ws.Range("D47").Value = "'0.6825"
$ws.Range("D48").Value = "'0.02661"
$ws.Range("D49").Value = "'0.00002108"
$ws.Range("D50").Value = "'0.01014"

# Column E ("Volume(1h)") text updates
$ws.Range("E17").Value = "16OneONE"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
